$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting old D:K (quarterly/annual data columns) to E:L.
# This mirrors the author adding a newer reporting period as the left-most data column.
$ws.Range("D:D").Insert()

# Copy formatting (number format, font, alignment) from the now-adjacent column E into
# the freshly inserted column D so every cell matches the style used by its row
# (date header row uses the date style, data rows use the numeric style).
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)

# Populate the shifted/refreshed data values for columns D through L across all rows.
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43100
$ws.Range("F7").Value = 42735
$ws.Range("G7").Value = 42369
$ws.Range("H7").Value = 42004
$ws.Range("I7").Value = 41639
$ws.Range("J7").Value = 41274
$ws.Range("K7").Value = 40908
$ws.Range("D8").Value = 27658200
$ws.Range("E8").Value = 27838800
$ws.Range("F8").Value = 24620900
$ws.Range("G8").Value = 25146000
$ws.Range("H8").Value = 23723400
$ws.Range("I8").Value = 23896100
$ws.Range("J8").Value = 23414800
$ws.Range("K8").Value = 22675700
$ws.Range("D9").Value = 14281800
$ws.Range("E9").Value = 14170700
$ws.Range("F9").Value = 12054700
$ws.Range("G9").Value = 12579800
$ws.Range("H9").Value = 12404700
$ws.Range("I9").Value = 12316100
$ws.Range("J9").Value = 11678800
$ws.Range("K9").Value = 11199300
$ws.Range("D10").Value = 13376400
$ws.Range("E10").Value = 13668100
$ws.Range("F10").Value = 12566300
$ws.Range("G10").Value = 12566300
$ws.Range("H10").Value = 11318600
$ws.Range("I10").Value = 11580100
$ws.Range("J10").Value = 11736000
$ws.Range("K10").Value = 11476300
$ws.Range("D12").Value = 375900
$ws.Range("E12").Value = 383700
$ws.Range("F12").Value = 373600
$ws.Range("G12").Value = 344500
$ws.Range("H12").Value = 305200
$ws.Range("I12").Value = 308500
$ws.Range("J12").Value = 288400
$ws.Range("K12").Value = 273500
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("D14").Value = 915500
$ws.Range("E14").Value = -319800
$ws.Range("F14").Value = 129000
$ws.Range("G14").Value = 812300
$ws.Range("H14").Value = 21300
$ws.Range("I14").Value = 21300
$ws.Range("J14").Value = 156000
$ws.Range("K14").Value = 140900
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("D17").Value = 24580600
$ws.Range("E17").Value = 23697600
$ws.Range("F17").Value = 21356000
$ws.Range("G17").Value = 22670900
$ws.Range("H17").Value = 21310000
$ws.Range("I17").Value = 21508500
$ws.Range("J17").Value = 20332700
$ws.Range("K17").Value = 19472300
$ws.Range("D18").Value = 3077600
$ws.Range("E18").Value = 4141300
$ws.Range("F18").Value = 3265000
$ws.Range("G18").Value = 2475100
$ws.Range("H18").Value = 2413400
$ws.Range("I18").Value = 2387600
$ws.Range("J18").Value = 3082100
$ws.Range("K18").Value = 3203300
$ws.Range("D20").Value = 48200
$ws.Range("E20").Value = 16800
$ws.Range("F20").Value = -4500
$ws.Range("G20").Value = -7900
$ws.Range("H20").Value = -42600
$ws.Range("I20").Value = 6700
$ws.Range("J20").Value = -64000
$ws.Range("K20").Value = -62200
$ws.Range("D21").Value = 4926500
$ws.Range("E21").Value = 5253500
$ws.Range("F21").Value = 4144500
$ws.Range("G21").Value = 3836000
$ws.Range("H21").Value = 3446000
$ws.Range("I21").Value = 3192900
$ws.Range("J21").Value = 3771700
$ws.Range("K21").Value = 3889300
$ws.Range("D22").Value = 440900
$ws.Range("E22").Value = 464500
$ws.Range("F22").Value = 309700
$ws.Range("G22").Value = 307400
$ws.Range("H22").Value = 307400
$ws.Range("I22").Value = 301800
$ws.Range("J22").Value = 274900
$ws.Range("K22").Value = 282900
$ws.Range("D23").Value = 2684900
$ws.Range("E23").Value = 3693600
$ws.Range("F23").Value = 2950800
$ws.Range("G23").Value = 2159800
$ws.Range("H23").Value = 2063300
$ws.Range("I23").Value = 2092500
$ws.Range("J23").Value = 2743300
$ws.Range("K23").Value = 2858200
$ws.Range("D24").Value = 803300
$ws.Range("E24").Value = 944700
$ws.Range("F24").Value = 902100
$ws.Range("G24").Value = 702400
$ws.Range("H24").Value = 672100
$ws.Range("I24").Value = 677700
$ws.Range("J24").Value = 798900
$ws.Range("K24").Value = 734800
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("D26").Value = 1881600
$ws.Range("E26").Value = 2748900
$ws.Range("F26").Value = 2048800
$ws.Range("G26").Value = 1457500
$ws.Range("H26").Value = 1391300
$ws.Range("I26").Value = 1414800
$ws.Range("J26").Value = 1944400
$ws.Range("K26").Value = 2123400
$ws.Range("D27").Value = 2619800
$ws.Range("E27").Value = 2745500
$ws.Range("F27").Value = 1929800
$ws.Range("G27").Value = 1438400
$ws.Range("H27").Value = 1255500
$ws.Range("I27").Value = 1595500
$ws.Range("J27").Value = 1876000
$ws.Range("K27").Value = 1961400
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 0
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("D32").Value = -48200
$ws.Range("E32").Value = -16800
$ws.Range("F32").Value = 4500
$ws.Range("G32").Value = 7900
$ws.Range("H32").Value = 42600
$ws.Range("I32").Value = -6700
$ws.Range("J32").Value = 64000
$ws.Range("K32").Value = 62200
$ws.Range("D33").Value = 2619800
$ws.Range("E33").Value = 2745500
$ws.Range("F33").Value = 1929800
$ws.Range("G33").Value = 1438400
$ws.Range("H33").Value = 1255500
$ws.Range("I33").Value = 1595500
$ws.Range("J33").Value = 1876000
$ws.Range("K33").Value = 1961400
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("D35").Value = 2619800
$ws.Range("E35").Value = 2745500
$ws.Range("F35").Value = 1929800
$ws.Range("G35").Value = 1438400
$ws.Range("H35").Value = 1255500
$ws.Range("I35").Value = 1595500
$ws.Range("J35").Value = 1876000
$ws.Range("K35").Value = 1961400
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43100
$ws.Range("F38").Value = 42735
$ws.Range("G38").Value = 42369
$ws.Range("H38").Value = 42004
$ws.Range("I38").Value = 41639
$ws.Range("J38").Value = 41274
$ws.Range("K38").Value = 40908
$ws.Range("D41").Value = 941300
$ws.Range("E41").Value = 715800
$ws.Range("F41").Value = 624900
$ws.Range("G41").Value = 582300
$ws.Range("H41").Value = 987400
$ws.Range("I41").Value = 1087200
$ws.Range("J41").Value = 1423800
$ws.Range("K41").Value = 1205500
$ws.Range("D42").Value = 4711200
$ws.Range("E42").Value = 3884300
$ws.Range("F42").Value = 14656600
$ws.Range("G42").Value = 2820700
$ws.Range("H42").Value = 2599700
$ws.Range("I42").Value = 3211100
$ws.Range("J42").Value = 1961200
$ws.Range("K42").Value = 1307600
$ws.Range("D43").Value = 3728400
$ws.Range("E43").Value = 3943800
$ws.Range("F43").Value = 3631900
$ws.Range("G43").Value = 3489400
$ws.Range("H43").Value = 2165400
$ws.Range("I43").Value = 2055500
$ws.Range("J43").Value = 4365700
$ws.Range("K43").Value = 2372300
$ws.Range("D44").Value = 2007200
$ws.Range("E44").Value = 1871500
$ws.Range("F44").Value = 1548300
$ws.Range("G44").Value = 1541600
$ws.Range("H44").Value = 1503500
$ws.Range("I44").Value = 1404700
$ws.Range("J44").Value = 1282400
$ws.Range("K44").Value = 1244200
$ws.Range("D45").Value = 206400
$ws.Range("E45").Value = 391600
$ws.Range("F45").Value = 982900
$ws.Range("G45").Value = 684400
$ws.Range("H45").Value = 1132100
$ws.Range("I45").Value = 1049100
$ws.Range("J45").Value = 1834500
$ws.Range("K45").Value = 1044700
$ws.Range("D46").Value = 11594600
$ws.Range("E46").Value = 10807000
$ws.Range("F46").Value = 21444600
$ws.Range("G46").Value = 8973700
$ws.Range("H46").Value = 8388000
$ws.Range("I46").Value = 8807600
$ws.Range("J46").Value = 7767500
$ws.Range("K46").Value = 7174300
$ws.Range("D47").Value = 2672600
$ws.Range("E47").Value = 3296400
$ws.Range("F47").Value = 3386200
$ws.Range("G47").Value = 3541000
$ws.Range("H47").Value = 2766800
$ws.Range("I47").Value = 1791800
$ws.Range("J47").Value = 1488900
$ws.Range("K47").Value = 1610500
$ws.Range("D48").Value = 6928300
$ws.Range("E48").Value = 6737500
$ws.Range("F48").Value = 5650300
$ws.Range("G48").Value = 5331700
$ws.Range("H48").Value = 7181900
$ws.Range("I48").Value = 4862700
$ws.Range("J48").Value = "NA"
$ws.Range("K48").Value = 4596600
$ws.Range("D49").Value = 27427000
$ws.Range("E49").Value = 27988000
$ws.Range("F49").Value = 17730800
$ws.Range("G49").Value = 17703900
$ws.Range("H49").Value = 18214400
$ws.Range("I49").Value = 18297400
$ws.Range("J49").Value = 36093300
$ws.Range("K49").Value = 18872500
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("F50").Value = 0
$ws.Range("G50").Value = 0
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("D52").Value = 943600
$ws.Range("E52").Value = 829200
$ws.Range("F52").Value = 1098400
$ws.Range("G52").Value = 1152300
$ws.Range("H52").Value = 1109600
$ws.Range("I52").Value = 941300
$ws.Range("J52").Value = 1256600
$ws.Range("K52").Value = 1112800
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("D54").Value = 49566200
$ws.Range("E54").Value = 49658200
$ws.Range("F54").Value = 49310300
$ws.Range("G54").Value = 36702500
$ws.Range("H54").Value = 35619800
$ws.Range("I54").Value = 34700900
$ws.Range("J54").Value = 33140200
$ws.Range("K54").Value = 33366700
$ws.Range("D57").Value = 4123300
$ws.Range("E57").Value = 4380200
$ws.Range("F57").Value = 4232100
$ws.Range("G57").Value = 4066100
$ws.Range("H57").Value = 3740700
$ws.Range("I57").Value = 3644200
$ws.Range("J57").Value = 6514300
$ws.Range("K57").Value = 3091800
$ws.Range("D58").Value = 3388400
$ws.Range("E58").Value = 3616200
$ws.Range("F58").Value = 2376400
$ws.Range("G58").Value = 2663600
$ws.Range("H58").Value = 2616500
$ws.Range("I58").Value = 2350600
$ws.Range("J58").Value = 2079000
$ws.Range("K58").Value = 2273700
$ws.Range("D59").Value = 3668900
$ws.Range("E59").Value = 3739600
$ws.Range("F59").Value = 3545500
$ws.Range("G59").Value = 3815900
$ws.Range("H59").Value = 5563900
$ws.Range("I59").Value = 5906200
$ws.Range("J59").Value = 6644400
$ws.Range("K59").Value = 2806600
$ws.Range("D60").Value = 11180600
$ws.Range("E60").Value = 11734900
$ws.Range("F60").Value = 10154000
$ws.Range("G60").Value = 10324600
$ws.Range("H60").Value = 11921100
$ws.Range("I60").Value = 11900900
$ws.Range("J60").Value = 9596400
$ws.Range("K60").Value = 8172100
$ws.Range("D61").Value = 16017500
$ws.Range("E61").Value = 17423400
$ws.Range("F61").Value = 20687300
$ws.Range("G61").Value = 8790800
$ws.Range("H61").Value = 6999000
$ws.Range("I61").Value = 7378200
$ws.Range("J61").Value = 4983900
$ws.Range("K61").Value = 4027300
$ws.Range("D62").Value = 3883200
$ws.Range("E62").Value = 5325000
$ws.Range("F62").Value = 3665500
$ws.Range("G62").Value = 3372700
$ws.Range("H62").Value = 3521900
$ws.Range("I62").Value = 3383900
$ws.Range("J62").Value = 5455100
$ws.Range("K62").Value = 6849200
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("F63").Value = 0
$ws.Range("G63").Value = 0
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 0
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 0
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("D66").Value = 31228300
$ws.Range("E66").Value = 33438700
$ws.Range("F66").Value = 34602200
$ws.Range("G66").Value = 22558700
$ws.Range("H66").Value = 22497000
$ws.Range("I66").Value = 22702300
$ws.Range("J66").Value = 19462000
$ws.Range("K66").Value = 19163600
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 0
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 0
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("F70").Value = 0
$ws.Range("G70").Value = 0
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("F71").Value = 0
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("D72").Value = 17835200
$ws.Range("E72").Value = 16467400
$ws.Range("F72").Value = 13503100
$ws.Range("G72").Value = 12851300
$ws.Range("H72").Value = 13258600
$ws.Range("I72").Value = 12513600
$ws.Range("J72").Value = 12258900
$ws.Range("K72").Value = 11963500
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("F73").Value = 0
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("F74").Value = 0
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("F75").Value = 0
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("D76").Value = 18337800
$ws.Range("E76").Value = 16219500
$ws.Range("F76").Value = 14708200
$ws.Range("G76").Value = 14143800
$ws.Range("H76").Value = 13122800
$ws.Range("I76").Value = 11998600
$ws.Range("J76").Value = 13678200
$ws.Range("K76").Value = 14203100
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43100
$ws.Range("F80").Value = 42735
$ws.Range("G80").Value = 42369
$ws.Range("H80").Value = 42004
$ws.Range("I80").Value = 41639
$ws.Range("J80").Value = 41274
$ws.Range("K80").Value = 40908
$ws.Range("D81").Value = 2619800
$ws.Range("E81").Value = 2745500
$ws.Range("F81").Value = 1929800
$ws.Range("G81").Value = 1438400
$ws.Range("H81").Value = 1255500
$ws.Range("I81").Value = 1595500
$ws.Range("J81").Value = 1876000
$ws.Range("K81").Value = 1961400
$ws.Range("D83").Value = 1796300
$ws.Range("E83").Value = 1092800
$ws.Range("F83").Value = 881900
$ws.Range("G83").Value = 1365500
$ws.Range("H83").Value = 1072600
$ws.Range("I83").Value = 796600
$ws.Range("J83").Value = 751700
$ws.Range("K83").Value = 747700
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("F84").Value = 0
$ws.Range("G84").Value = 0
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("F85").Value = 0
$ws.Range("G85").Value = 0
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("F86").Value = 0
$ws.Range("G86").Value = 0
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("F87").Value = 0
$ws.Range("G87").Value = 0
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("F88").Value = 0
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("D89").Value = 3490500
$ws.Range("E89").Value = 3317700
$ws.Range("F89").Value = 2976600
$ws.Range("G89").Value = 2656900
$ws.Range("H89").Value = 2456000
$ws.Range("I89").Value = 2643400
$ws.Range("J89").Value = 3206600
$ws.Range("K89").Value = 3057800
$ws.Range("D91").Value = -1055800
$ws.Range("E91").Value = -1087200
$ws.Range("F91").Value = -1037800
$ws.Range("G91").Value = -1051300
$ws.Range("H91").Value = -1104000
$ws.Range("I91").Value = -1165700
$ws.Range("J91").Value = -1095100
$ws.Range("K91").Value = -1038800
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("F92").Value = 0
$ws.Range("G92").Value = 0
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("F93").Value = 0
$ws.Range("G93").Value = 0
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("D94").Value = 365800
$ws.Range("E94").Value = -12832200
$ws.Range("F94").Value = -951400
$ws.Range("G94").Value = -1704300
$ws.Range("H94").Value = -2205800
$ws.Range("I94").Value = -2419000
$ws.Range("J94").Value = -1201700
$ws.Range("K94").Value = -900300
$ws.Range("D96").Value = -572200
$ws.Range("E96").Value = -409500
$ws.Range("F96").Value = -1210600
$ws.Range("G96").Value = -461100
$ws.Range("H96").Value = -467900
$ws.Range("I96").Value = -951400
$ws.Range("J96").Value = -936900
$ws.Range("K96").Value = -919100
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("F97").Value = 0
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("F99").Value = 0
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("D100").Value = -3647600
$ws.Range("E100").Value = 9300200
$ws.Range("F100").Value = -1813100
$ws.Range("G100").Value = -1101800
$ws.Range("H100").Value = -400600
$ws.Range("I100").Value = -384800
$ws.Range("J100").Value = -1665000
$ws.Range("K100").Value = -2188000
$ws.Range("D101").Value = 15700
$ws.Range("E101").Value = 305200
$ws.Range("F101").Value = -169400
$ws.Range("G101").Value = -255800
$ws.Range("H101").Value = 50500
$ws.Range("I101").Value = -176200
$ws.Range("J101").Value = -68400
$ws.Range("K101").Value = -1200
$ws.Range("D102").Value = 224400
$ws.Range("E102").Value = 90900
$ws.Range("F102").Value = 42600
$ws.Range("G102").Value = -405000
$ws.Range("H102").Value = -99900
$ws.Range("I102").Value = -336600
$ws.Range("J102").Value = 271500
$ws.Range("K102").Value = -31700

